$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 363. This shifts all existing rows
# 363..494 down to 364..495 (and the sheet dimension grows to R495),
# exactly matching the diff (old row 363 -> new row 364, ..., old
# row 494 -> new row 495).
$ws.Rows.Item(363).Insert()

# Populate the newly inserted row 363 with its data. Columns
# A, B, C, E, F, G, H, I, N, O, Q, R keep the same values the block
# already used; D, J, K, L, M, P carry the new figures from the diff.
$ws.Range("A363").Value = 8
$ws.Range("B363").Value = "Terminal La Palmera de La Serena"
$ws.Range("C363").Value = "Coquimbo"
$ws.Range("D363").Value = 45229
$ws.Range("D363").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E363").Value = 4
$ws.Range("F363").Value = 100112012
$ws.Range("G363").Value = "Espinaca"
$ws.Range("H363").Value = "Sin especificar"
$ws.Range("I363").Value = "Primera"
$ws.Range("J363").Value = 1200
$ws.Range("K363").Value = 450
$ws.Range("L363").Value = 500
$ws.Range("M363").Value = 475
$ws.Range("N363").Value = "`$/atado 300 a 500 gramos"
$ws.Range("O363").Value = "Provincia del Elquí"
$ws.Range("P363").Value = 950
$ws.Range("Q363").Value = 0.5
$ws.Range("R363").Value = "Hortaliza"
